# feat: add 2022-Q4 data
#
# Inserts a new "2022-Q4" sheet (copied from the existing "2022-Q3" sheet to
# keep the same layout/formatting), fills it with the new quarter's fund
# data, and updates the "总计" (Total) summary sheet with a new leading row
# for 2022-Q4 (shifting the existing quarters down).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Duplicate the "2022-Q3" sheet, inserting the copy right before it, then
#    rename the copy to "2022-Q4". Copying (rather than Worksheets.Add)
#    keeps every cell's existing style/border formatting intact.
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3, $null)

$q4 = $wb.Worksheets.Item("2022-Q3 (2)")
$q4.Name = "2022-Q4"

# ---------------------------------------------------------------------------
# 2. Overwrite the new sheet's data rows (2-8) with the 2022-Q4 figures.
#    Row/column 1 (headers) and column A (row index) are unchanged from the
#    copy, so only columns B-H for rows 2-8 need to be written.
#
#    Columns B and D-G hold digit-string values (fund codes / percentages)
#    that must stay text (e.g. the leading zero in "010695", or the trailing
#    zero in "0.7710") instead of being auto-converted to numbers, so they
#    are entered with a leading apostrophe - the normal Excel way to force
#    text entry - which is stripped from the stored value automatically.
# ---------------------------------------------------------------------------
$q4.Range("B2").Value = "'010695"
$q4.Range("C2").Value = "华夏磐益一年定期开放混合"
$q4.Range("D2").Value = "'16.03"
$q4.Range("E2").Value = "'98.69"
$q4.Range("F2").Value = "'4.81"
$q4.Range("G2").Value = "'0.7710"
$q4.Range("H2").Value = 4

$q4.Range("B3").Value = "'501079"
$q4.Range("C3").Value = "大成科创主题混合（LOF）A"
$q4.Range("D3").Value = "'9.55"
$q4.Range("E3").Value = "'80.68"
$q4.Range("F3").Value = "'6.20"
$q4.Range("G3").Value = "'0.5921"
$q4.Range("H3").Value = 5

$q4.Range("B4").Value = "'012473"
$q4.Range("C4").Value = "大成成长回报六个月持有混合A"
$q4.Range("D4").Value = "'7.49"
$q4.Range("E4").Value = "'75.28"
$q4.Range("F4").Value = "'5.91"
$q4.Range("G4").Value = "'0.4427"
$q4.Range("H4").Value = 4

$q4.Range("B5").Value = "'010371"
$q4.Range("C5").Value = "大成成长进取混合A"
$q4.Range("D5").Value = "'3.61"
$q4.Range("E5").Value = "'74.21"
$q4.Range("F5").Value = "'5.15"
$q4.Range("G5").Value = "'0.1859"
$q4.Range("H5").Value = 5

$q4.Range("B6").Value = "'010372"
$q4.Range("C6").Value = "大成成长进取混合C"
$q4.Range("D6").Value = "'1.50"
$q4.Range("E6").Value = "'74.21"
$q4.Range("F6").Value = "'5.15"
$q4.Range("G6").Value = "'0.0772"
$q4.Range("H6").Value = 5

$q4.Range("B7").Value = "'012474"
$q4.Range("C7").Value = "大成成长回报六个月持有混合C"
$q4.Range("D7").Value = "'0.37"
$q4.Range("E7").Value = "'75.28"
$q4.Range("F7").Value = "'5.91"
$q4.Range("G7").Value = "'0.0219"
$q4.Range("H7").Value = 4

$q4.Range("B8").Value = "'016198"
$q4.Range("C8").Value = "大成科创主题混合（LOF）C"
$q4.Range("D8").Value = "'0.06"
$q4.Range("E8").Value = "'80.68"
$q4.Range("F8").Value = "'6.20"
$q4.Range("G8").Value = "'0.0037"
$q4.Range("H8").Value = 5

# ---------------------------------------------------------------------------
# 3. Update the "总计" summary sheet: shift the existing quarter rows down by
#    one and insert the new 2022-Q4 row at the top of the data (row 2). The
#    sheet is small, so rows are rewritten in place (bottom-up) rather than
#    using Rows.Insert, which would otherwise drag stray border formatting
#    into the new row.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Grab a copy of row 5's formatting before overwriting it, so the brand new
# row 6 (2021-Q3) can get the same column-A style ("s=2") the other data
# rows have.
$total.Range("A5").Copy()
$total.Range("A6").PasteSpecial(-4122)

$total.Range("A6").Value = 4
$total.Range("B6").Value = "2021-Q3"
$total.Range("C6").Value = 1
$total.Range("D6").Value = 0.26

$total.Range("A5").Value = 3
$total.Range("B5").Value = "2022-Q1"
$total.Range("C5").Value = 5
$total.Range("D5").Value = 1.65

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2022-Q2"
$total.Range("C4").Value = 8
$total.Range("D4").Value = 1.37

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q3"
$total.Range("C3").Value = 7
$total.Range("D3").Value = 2.6

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 7
$total.Range("D2").Value = 2.09

# ---------------------------------------------------------------------------
# 4. Restore the original active/selected sheet ("2021-Q3", the last tab)
#    since copying/renaming sheets above shifts Excel's active-sheet focus
#    to the newly inserted sheet.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q3").Activate()
